# "need to close the year" - update progress/status figures and comments
# on the Sheet1 deliverables tracker for the end-of-year close out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: the mWater-account-naming task is now fully done; drop the
# "2 accounts from 48..." outstanding-work comment.
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = ""

# Row 8: HANWASH portal structure work is nearly done (99%) and is an
# ongoing task - flag it in bold.
$ws.Range("G8").Value = 0.99
$ws.Range("H8").Value = "Ongoing task"
$ws.Range("H8").Font.Bold = $true

# Row 10: user guides are fully completed now; update comment accordingly.
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = "Guides completed and shared"

# Row 11: stakeholder account creation effort wrapped up; grow the row to
# fit the longer closing comment.
$ws.Rows.Item(11).RowHeight = 119
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = "Elements were shared and invitation were shared. We will continue with to relaunch to process for the rest of the non responsive members"

# Move the active selection to reflect where the author left off reviewing.
$ws.Range("G7").Select()
